$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2282.3462
$ws.Range("J40").Value = 2347.5417
$ws.Range("L40").Value = 2347.5417
$ws.Range("N40").Value = -2697.5417

# Row 96
$ws.Range("H96").Value = 3611.1667
$ws.Range("I96").Value = 3500.1538
$ws.Range("J96").Value = 3899.8
$ws.Range("K96").Value = 10500.4614
$ws.Range("L96").Value = 11699.4
$ws.Range("M96").Value = -9127.4614
$ws.Range("N96").Value = -14445.4

# Row 132
$ws.Range("H132").Value = 296319.66
$ws.Range("I132").Value = 2167.7188
$ws.Range("K132").Value = 6503.1564
$ws.Range("M132").Value = -3973.1564

# Row 141
$ws.Range("H141").Value = 3935.8
$ws.Range("I141").Value = 2419.75
$ws.Range("K141").Value = 7259.25
$ws.Range("M141").Value = -2079.25

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3600
$ws.Range("I61").Value = 3900
$ws.Range("K61").Value = 3900
$ws.Range("M61").Value = -3688

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").ClearContents()

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").ClearContents()

# Row 80
$ws.Range("H80").Value = 19904.54
$ws.Range("J80").Value = 19904.54
$ws.Range("L80").Value = 19904.54
$ws.Range("N80").Value = -21900.54

# Row 83
$ws.Range("H83").Value = 19904.54
$ws.Range("J83").Value = 19904.54
$ws.Range("L83").Value = 59713.62
$ws.Range("N83").Value = -69697.62

# Row 88
$ws.Range("H88").Value = 2970.7
$ws.Range("I88").Value = 2250
$ws.Range("J88").Value = 3150.875
$ws.Range("K88").Value = 2250
$ws.Range("L88").Value = 3150.875
$ws.Range("M88").Value = -1844
$ws.Range("N88").Value = -3962.875

# Row 91
$ws.Range("H91").Value = 2970.7
$ws.Range("I91").Value = 2250
$ws.Range("J91").Value = 3150.875
$ws.Range("K91").Value = 2250
$ws.Range("L91").Value = 3150.875
$ws.Range("M91").Value = -846
$ws.Range("N91").Value = -5958.875

# Row 102
$ws.Range("H102").Value = 4674.3076
$ws.Range("I102").Value = 4750.5557
$ws.Range("J102").Value = 4502.75
$ws.Range("K102").Value = 4750.5557
$ws.Range("L102").Value = 4502.75
$ws.Range("M102").Value = -3128.5557
$ws.Range("N102").Value = -7746.75

# Row 136
$ws.Range("H136").Value = 3600
$ws.Range("I136").Value = 3900
$ws.Range("K136").Value = 11700
$ws.Range("M136").Value = -9150

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 426.33334
$ws.Range("J22").Value = 499.5
$ws.Range("L22").Value = 499.5
$ws.Range("N22").Value = -845.5

# Row 55
$ws.Range("H55").Value = 40000
$ws.Range("J55").Value = 40000
$ws.Range("L55").Value = 40000
$ws.Range("N55").Value = -40546

# Row 86
$ws.Range("H86").Value = 2807.2727
$ws.Range("I86").Value = 2320
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2320
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1197
$ws.Range("N86").Value = -7246

# Row 87
$ws.Range("H87").Value = 37677
$ws.Range("J87").Value = 37677
$ws.Range("L87").Value = 37677
$ws.Range("N87").Value = -40173

# Row 89
$ws.Range("H89").Value = 2807.2727
$ws.Range("I89").Value = 2320
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 11600
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -5984
$ws.Range("N89").Value = -36232

# Row 90
$ws.Range("H90").Value = 37677
$ws.Range("J90").Value = 37677
$ws.Range("L90").Value = 113031
$ws.Range("N90").Value = -125511

# Row 134
$ws.Range("H134").Value = 98438.09
$ws.Range("I134").Value = 6583.3335
$ws.Range("J134").Value = 135180
$ws.Range("K134").Value = 19750.0005
$ws.Range("L134").Value = 405540
$ws.Range("M134").Value = -17215.0005
$ws.Range("N134").Value = -410610

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 9615.385
$ws.Range("I62").Value = 10312.5
$ws.Range("J62").Value = 8500
$ws.Range("K62").Value = 10312.5
$ws.Range("L62").Value = 8500
$ws.Range("M62").Value = -9688.5
$ws.Range("N62").Value = -9748

# Row 65
$ws.Range("H65").Value = 9615.385
$ws.Range("I65").Value = 10312.5
$ws.Range("J65").Value = 8500
$ws.Range("K65").Value = 51562.5
$ws.Range("L65").Value = 42500
$ws.Range("M65").Value = -48442.5
$ws.Range("N65").Value = -48740

# Row 122
$ws.Range("H122").Value = 1675168.5
$ws.Range("J122").Value = 12499.75
$ws.Range("L122").Value = 37499.25
$ws.Range("N122").Value = -42399.25

# Row 130
$ws.Range("H130").Value = 36992
$ws.Range("J130").Value = 36992
$ws.Range("L130").Value = 36992
$ws.Range("N130").Value = -47032

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 4902162
$ws.Range("I2").Value = 14705900
$ws.Range("J2").Value = 293.25
$ws.Range("K2").Value = 88235400
$ws.Range("L2").Value = 1759.5
$ws.Range("M2").Value = -88235287
$ws.Range("N2").Value = -1985.5

# Row 12
$ws.Range("I12").Value = 27.75
$ws.Range("J12").Value = 30.083334
$ws.Range("K12").Value = 83.25
$ws.Range("L12").Value = 90.25000199999999
$ws.Range("M12").Value = 89.75
$ws.Range("N12").Value = -436.250002

# Row 38
$ws.Range("H38").Value = 289.33334
$ws.Range("I38").Value = 80.40000000000001
$ws.Range("J38").Value = 550.5
$ws.Range("K38").Value = 241.2
$ws.Range("L38").Value = 1651.5
$ws.Range("M38").Value = 105.8
$ws.Range("N38").Value = -2345.5

# Row 110
$ws.Range("H110").Value = 9122.700000000001
$ws.Range("I110").Value = 306.75
$ws.Range("K110").Value = 920.25
$ws.Range("M110").Value = 3169.75

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1951.5
$ws.Range("I7").Value = 2042.4
$ws.Range("J7").Value = 1800
$ws.Range("K7").Value = 2042.4
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = -1930.4
$ws.Range("N7").Value = -2024

# Row 16
$ws.Range("H16").Value = 2239.8

# Row 40
$ws.Range("H40").Value = 3049.5
$ws.Range("I40").Value = 3999
$ws.Range("J40").Value = 2100
$ws.Range("K40").Value = 3999
$ws.Range("L40").Value = 2100
$ws.Range("M40").Value = -3863
$ws.Range("N40").Value = -2372

# Row 50
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("N50").ClearContents()

# Row 93
$ws.Range("H93").Value = 1284.6538
$ws.Range("I93").Value = 968.4737
$ws.Range("J93").Value = 2142.8572
$ws.Range("K93").Value = 968.4737
$ws.Range("L93").Value = 2142.8572
$ws.Range("M93").Value = 279.5263
$ws.Range("N93").Value = -4638.8572

# Row 122
$ws.Range("H122").Value = 8601
$ws.Range("I122").Value = 9032.23
$ws.Range("J122").Value = 2995
$ws.Range("K122").Value = 27096.69
$ws.Range("L122").Value = 8985
$ws.Range("M122").Value = -24646.69
$ws.Range("N122").Value = -13885

# Row 126
$ws.Range("H126").Value = 1951.5
$ws.Range("I126").Value = 2042.4
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 6127.200000000001
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -3657.200000000001
$ws.Range("N126").Value = -10340
